# fix: update SNAPP_TOKEN in .env file
# Refresh scraped price/checksum data for the product inventory sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2 - MOv6kw
$ws.Range("I2").Value = 12239471
$ws.Range("U2").Value = "c9db631f8c37"

# Row 3 - exLEv4
$ws.Range("I3").Value = 12450497
$ws.Range("J3").Value = 12450497
$ws.Range("U3").Value = "61428cd16d16"

# Row 4 - za254K
$ws.Range("I4").Value = 13169606
$ws.Range("J4").Value = 13169606
$ws.Range("U4").Value = "a6a714f70d53"

# Row 5 - a8bOMv
$ws.Range("I5").Value = 13408227
$ws.Range("J5").Value = 13408227
$ws.Range("U5").Value = "d05f058d66d3"

# Row 6 - Z4bQR3
$ws.Range("I6").Value = 11674573
$ws.Range("J6").Value = 11674573
$ws.Range("U6").Value = "164a286afaa7"

# Row 7 - b1byEJ
$ws.Range("I7").Value = 14270185
$ws.Range("J7").Value = 14270185
$ws.Range("U7").Value = "8c5312a5e111"

# Row 8 - dJbV8l (now wins buybox)
$ws.Range("I8").Value = 8614705
$ws.Range("J8").Value = 8614705
$ws.Range("K8").Value = "بله"
$ws.Range("L8").Value = 1
$ws.Range("U8").Value = "3c989061899b"

# Row 9 - X9brx7
$ws.Range("I9").Value = 10340243
$ws.Range("J9").Value = 10340243
$ws.Range("U9").Value = "dc37d0e519e4"
$ws.Range("W9").Value = 4950000
